$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows per re-pulled data / mean calculation
$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = -6
$ws.Range("F8").Value = -1
$ws.Range("F12").Value = -1
